# First round of optimising len/ac1thresh/correlThresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Make room for two new columns (L,M) before the existing "no of
# trades / profit / profit factor" block, which slides from L:N to N:P. ---
$ws.Range("L1:M1").EntireColumn.Insert()

# New column headers for the inserted L/M columns.
$ws.Range("L1").Value = "from"
$ws.Range("M1").Value = "processed bars"

# Brand new "note" column at the end (Q).
$ws.Range("Q1").Value = "note"

# --- Row 2: existing inputs get new values, plus a few new ones. ---
$ws.Range("B2").Value = "eur/usd"
$ws.Range("C2").Value = "M30"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 84
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 25
$ws.Range("I2").Value = "ac1"
$ws.Range("J2").Value = 6000

# "from" date, formatted as a date (built-in date format -> numFmtId 14).
$ws.Range("L2").NumberFormat = "mm-dd-yy"
$ws.Range("L2").Value = 42005

$ws.Range("M2").Value = 10000
$ws.Range("N2").Value = 60
$ws.Range("O2").Value = 1870
$ws.Range("P2").Value = 2.6

# --- Row 3: "best len" result row. ---
$ws.Range("A3").Value = "best len"
$ws.Range("D3").Value = 10
$ws.Range("L3").NumberFormat = "mm-dd-yy"

# --- Row 4: "best ac1 thresh" result row. ---
$ws.Range("A4").Value = "best ac1 thresh"
$ws.Range("G4").Value = 30
$ws.Range("N4").Value = 238
$ws.Range("O4").Value = 4360
$ws.Range("P4").Value = 1.66
$ws.Range("Q4").Value = "!no of trades drops dramaticaly with higher ac1, but with better profit factor"

# --- View state: scrolled right, with Q6 selected. ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("Q6").Select()
